$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.487.16"
$ws.Range("E2").Value = "  +3.28%  "
$ws.Range("D3").Value = "2.315.47"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'517.03"
$ws.Range("E5").Value = "  +2.79%  "
$ws.Range("D6").Value = "'135.32"
$ws.Range("E6").Value = "  +6.13%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("D9").Value = "2.335.24"
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("D10").Value = "'0.102"
$ws.Range("E10").Value = "  +4.96%  "
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("E12").Value = "  +4.78%  "
$ws.Range("D13").Value = "'0.341"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("D14").Value = "'23.94"
$ws.Range("E14").Value = "  +2.39%  "
$ws.Range("D15").Value = "2.728.84"
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("D16").Value = "56.600.61"
$ws.Range("E16").Value = "  +3.25%  "
$ws.Range("E17").Value = "  +3.04%  "
$ws.Range("D18").Value = "2.334.00"
$ws.Range("E18").Value = "  +2.62%  "
$ws.Range("D19").Value = "'10.48"
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("E20").Value = "  +1.24%  "
$ws.Range("D21").Value = "'322.86"
$ws.Range("E21").Value = "  +4.22%  "
$ws.Range("D22").Value = "'6.57"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "'0.996"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").Value = "'60.65"
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("D25").Value = "'0.164"
$ws.Range("E25").Value = "  +6.34%  "
$ws.Range("D26").Value = "'0.994"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").Value = "'7.98"
$ws.Range("E27").Value = "  +7.26%  "
$ws.Range("D28").Value = "'1.28"
$ws.Range("E28").Value = "  +12.78%  "
$ws.Range("D29").Value = "0.0₃0739"
$ws.Range("E29").Value = "  +5.75%  "
$ws.Range("E30").Value = "  +4.68%  "
$ws.Range("D31").Value = "'166.63"
$ws.Range("E31").Value = "  -3.20%  "
$ws.Range("D32").Value = "'6.20"
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("E33").Value = "  +2.68%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "'0.992"
$ws.Range("E35").Value = "  -0.46%  "
$ws.Range("E36").Value = "  +3.05%  "
$ws.Range("D37").Value = "'0.922"
$ws.Range("E37").Value = "  +2.57%  "
$ws.Range("D38").Value = "'4.01"
$ws.Range("E38").Value = "  +4.32%  "
$ws.Range("E39").Value = "  +7.37%  "
$ws.Range("D40").Value = "'37.94"
$ws.Range("E40").Value = "  +3.30%  "
$ws.Range("D41").Value = "'0.383"
$ws.Range("E41").Value = "  +2.48%  "
$ws.Range("D42").Value = "'140.07"
$ws.Range("E42").Value = "  +4.16%  "
$ws.Range("D43").Value = "'3.60"
$ws.Range("E43").Value = "  +4.44%  "
$ws.Range("D44").Value = "'5.24"
$ws.Range("E44").Value = "  +8.43%  "
$ws.Range("D45").Value = "'276.36"
$ws.Range("E45").Value = "  +7.68%  "
$ws.Range("D46").Value = "'0.0932"
$ws.Range("E46").Value = "  +2.35%  "
$ws.Range("D47").Value = "'0.0507"
$ws.Range("E47").Value = "  +0.89%  "
$ws.Range("D48").Value = "'0.560"
$ws.Range("E48").Value = "  +2.95%  "
$ws.Range("E49").Value = "  +3.31%  "
$ws.Range("D50").Value = "'0.381"
$ws.Range("E50").Value = "  +2.27%  "
$ws.Range("D51").Value = "'17.77"
$ws.Range("E51").Value = "  +9.33%  "
